$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole D2:E51 data range to Text format first so that
# numeric-looking strings (e.g. "0.999", "1.00") are not auto-converted
# to actual numbers when we set their .Value below.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '64.945.99'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '3.515.62'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '587.26'
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('D6').Value = '133.78'
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('D7').Value = '3.513.26'
$ws.Range('E7').Value = '  -1.12%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = '0.125'
$ws.Range('E10').Value = '  +2.12%  '
$ws.Range('D11').Value = '7.15'
$ws.Range('E11').Value = '  +3.23%  '
$ws.Range('D12').Value = '0.385'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D13').Value = '4.109.05'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').Value = '27.61'
$ws.Range('E14').Value = '  +2.98%  '
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').Value = '3.513.51'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '64.936.34'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('D20').Value = '14.29'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').Value = '5.66'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').Value = '390.83'
$ws.Range('E22').Value = '  +0.90%  '
$ws.Range('D23').Value = '0.575'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '74.91'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('D25').Value = '3.654.98'
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('E27').Value = '  -2.51%  '
$ws.Range('E28').Value = '  +8.21%  '
$ws.Range('D29').Value = '7.56'
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').Value = '8.28'
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('D33').Value = '3.521.52'
$ws.Range('E33').Value = '  -1.06%  '
$ws.Range('D34').Value = '24.11'
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  +2.21%  '
$ws.Range('D37').Value = '5.18'
$ws.Range('E37').Value = '  +4.76%  '
$ws.Range('D38').Value = '1.57'
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('D39').Value = '169.74'
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('D40').Value = '6.95'
$ws.Range('E40').Value = '  +1.22%  '
$ws.Range('D41').Value = '0.0805'
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').Value = '0.819'
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('D43').Value = '26.01'
$ws.Range('E43').Value = '  -2.44%  '
$ws.Range('D44').Value = '42.93'
$ws.Range('E45').Value = '  +3.90%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '4.42'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').Value = '1.65'
$ws.Range('E48').Value = '  +1.42%  '
$ws.Range('D49').Value = '2.489.07'
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('D50').Value = '6.85'
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('E51').Value = '  +3.14%  '

# Remove the temporary Text number format again so the cells end up
# with the same (default/no explicit style) formatting as before the edit.
$dataRange.ClearFormats()

